$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($ws, $addr, $val)
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextCell $ws "D2" "30.300.76"
Set-TextCell $ws "E2" "  -0.73%  "
Set-TextCell $ws "D3" "2.082.99"
Set-TextCell $ws "E3" "  +3.26%  "
Set-TextCell $ws "D4" "1.000"
Set-TextCell $ws "E4" "  -0.21%  "
Set-TextCell $ws "D5" "327.71"
Set-TextCell $ws "E5" "  +0.94%  "
Set-TextCell $ws "D6" "0.9998"
Set-TextCell $ws "D7" "0.5201"
Set-TextCell $ws "E7" "  +1.62%  "
Set-TextCell $ws "E8" "  +3.29%  "
Set-TextCell $ws "D9" "0.08819"
Set-TextCell $ws "E9" "  +0.66%  "
Set-TextCell $ws "D10" "46.11"
Set-TextCell $ws "E10" "  +6.03%  "
Set-TextCell $ws "D11" "1.161"
Set-TextCell $ws "E11" "  +2.11%  "
Set-TextCell $ws "D12" "24.57"
Set-TextCell $ws "E12" "  -0.45%  "
Set-TextCell $ws "D13" "2.084.45"
Set-TextCell $ws "E13" "  +3.33%  "
Set-TextCell $ws "D14" "6.704"
Set-TextCell $ws "E14" "  +1.48%  "
Set-TextCell $ws "D15" "7.684"
Set-TextCell $ws "E15" "  +2.64%  "
Set-TextCell $ws "D16" "95.26"
Set-TextCell $ws "E16" "  +0.80%  "
Set-TextCell $ws "E17" "  -0.05%  "
Set-TextCell $ws "D18" "0.00001121"
Set-TextCell $ws "E18" "  +0.50%  "
Set-TextCell $ws "D19" "0.06623"
Set-TextCell $ws "E19" "  +1.40%  "
Set-TextCell $ws "D20" "18.85"
Set-TextCell $ws "E20" "  -0.83%  "
Set-TextCell $ws "D21" "1.0000"
Set-TextCell $ws "E21" "  -0.11%  "
Set-TextCell $ws "D22" "6.330"
Set-TextCell $ws "E22" "  +1.63%  "
Set-TextCell $ws "D23" "30.341.96"
Set-TextCell $ws "E23" "  -0.78%  "
Set-TextCell $ws "D24" "12.33"
Set-TextCell $ws "E24" "  +3.89%  "
Set-TextCell $ws "D25" "2.288"
Set-TextCell $ws "E25" "  +2.60%  "
Set-TextCell $ws "D26" "2.324.62"
Set-TextCell $ws "E26" "  +3.25%  "
Set-TextCell $ws "D27" "22.32"
Set-TextCell $ws "E27" "  -0.45%  "
Set-TextCell $ws "D28" "2.601"
Set-TextCell $ws "E28" "  +6.90%  "
Set-TextCell $ws "D29" "162.07"
Set-TextCell $ws "E29" "  -0.66%  "
Set-TextCell $ws "D30" "131.03"
Set-TextCell $ws "E30" "  -0.40%  "
Set-TextCell $ws "D31" "1.191"
Set-TextCell $ws "E31" "  +4.27%  "
Set-TextCell $ws "D32" "0.1066"
Set-TextCell $ws "E32" "  +1.30%  "
Set-TextCell $ws "D33" "1.640"
Set-TextCell $ws "E33" "  +20.59%  "
Set-TextCell $ws "D34" "6.198"
Set-TextCell $ws "D35" "3.821"
Set-TextCell $ws "E35" "  -0.34%  "
Set-TextCell $ws "B36" "FraxShare"
Set-TextCell $ws "C36" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextCell $ws "D36" "9.876"
Set-TextCell $ws "E36" "  +8.49%  "
Set-TextCell $ws "B37" "VeChain"
Set-TextCell $ws "C37" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextCell $ws "D37" "0.02581"
Set-TextCell $ws "E37" "  +2.09%  "
Set-TextCell $ws "E38" "  +3.50%  "
Set-TextCell $ws "D39" "0.06677"
Set-TextCell $ws "E39" "  +0.33%  "
Set-TextCell $ws "D40" "5.452"
Set-TextCell $ws "E40" "  -0.88%  "
Set-TextCell $ws "D41" "0.2249"
Set-TextCell $ws "E41" "  +2.32%  "
Set-TextCell $ws "D42" "0.6835"
Set-TextCell $ws "E42" "  +2.28%  "
Set-TextCell $ws "D43" "1.244"
Set-TextCell $ws "E43" "  +0.96%  "
Set-TextCell $ws "E44" "  -0.12%  "
Set-TextCell $ws "B45" "EnergySwap"
Set-TextCell $ws "C45" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell $ws "D45" "14.00"
Set-TextCell $ws "E45" "  +2.23%  "
Set-TextCell $ws "B46" "Decentraland"
Set-TextCell $ws "C46" "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextCell $ws "D46" "0.6358"
Set-TextCell $ws "E46" "  +2.62%  "
Set-TextCell $ws "D47" "2.205"
Set-TextCell $ws "E47" "  +0.20%  "
Set-TextCell $ws "D48" "3.608"
Set-TextCell $ws "E48" "  -1.55%  "
Set-TextCell $ws "D49" "1.248"
Set-TextCell $ws "E49" "  -1.73%  "
Set-TextCell $ws "D50" "1.188"
Set-TextCell $ws "E50" "  +7.07%  "
Set-TextCell $ws "D51" "81.68"
Set-TextCell $ws "E51" "  +0.60%  "
